$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row4 -> F4 871 -> 876, row6 -> F6 37 -> 39
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 876
$wsExhibit.Range("F6").Value = 39

# Sheet "全部类型" (sheet4): row5 -> F5 871 -> 876, row7 -> F7 37 -> 39
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 876
$wsAll.Range("F7").Value = 39
